$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The four data rows (2, 5, 6, 7) get their values cyclically rotated:
#   row2 <- old row7, row5 <- old row6, row6 <- old row2, row7 <- old row5
# Apply the resulting target values directly (date serial, volumes, prices,
# unit text, origin text, price/kg, kg/unit).

$ws.Range("D2").Value = 44698
$ws.Range("M2").Value = 120
$ws.Range("N2").Value = 16000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 16500
$ws.Range("Q2").Value = '$/caja 18 kilos granel'
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 917
$ws.Range("T2").Value = 18

$ws.Range("D5").Value = 44344
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 13000
$ws.Range("O5").Value = 14000
$ws.Range("P5").Value = 13500
$ws.Range("Q5").Value = '$/caja 18 kilos granel'
$ws.Range("R5").Value = "Provincia de Curicó"
$ws.Range("S5").Value = 750
$ws.Range("T5").Value = 18

$ws.Range("D6").Value = 44330
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 16000
$ws.Range("P6").Value = 15500
$ws.Range("Q6").Value = '$/caja 18 kilos granel'
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 861
$ws.Range("T6").Value = 18

$ws.Range("D7").Value = 44334
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 12500
$ws.Range("Q7").Value = '$/caja 12 kilos empedrada'
$ws.Range("R7").Value = "Región de O'Higgins"
$ws.Range("S7").Value = 1042
$ws.Range("T7").Value = 12
